{"js": "// Update the date heading paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst headingRange = paragraphs.items[0].getRange();\nheadingRange.insertText(\"2026-01-30 Friday\", Word.InsertLocation.replace);\n\n// Update the division problems in the table (row, col are 0-based).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  { row: 0, col: 0, text: \"23\u00f75=\" },\n  { row: 0, col: 1, text: \"21\u00f75=\" },\n  { row: 0, col: 2, text: \"27\u00f79=\" },\n  { row: 0, col: 3, text: \"37\u00f78=\" },\n  { row: 0, col: 4, text: \"30\u00f76=\" },\n\n  { row: 4, col: 0, text: \"64\u00f78=\" },\n  { row: 4, col: 1, text: \"22\u00f75=\" },\n  { row: 4, col: 2, text: \"44\u00f73=\" },\n  { row: 4, col: 3, text: \"51\u00f76=\" },\n  { row: 4, col: 4, text: \"89\u00f78=\" },\n\n  { row: 8, col: 0, text: \"75\u00f72=\" },\n  { row: 8, col: 1, text: \"47\u00f73=\" },\n  { row: 8, col: 2, text: \"72\u00f78=\" },\n  { row: 8, col: 3, text: \"42\u00f75=\" },\n  { row: 8, col: 4, text: \"70\u00f75=\" },\n\n  { row: 12, col: 0, text: \"19\u00f76=\" },\n  { row: 12, col: 1, text: \"97\u00f75=\" },\n  { row: 12, col: 2, text: \"61\u00f73=\" },\n  { row: 12, col: 3, text: \"70\u00f79=\" },\n  { row: 12, col: 4, text: \"83\u00f77=\" },\n\n  { row: 16, col: 0, text: \"70\u00f77=\" },\n  { row: 16, col: 1, text: \"99\u00f75=\" },\n  { row: 16, col: 2, text: \"61\u00f77=\" },\n  { row: 16, col: 3, text: \"66\u00f79=\" },\n  { row: 16, col: 4, text: \"99\u00f74=\" },\n];\n\nconst cellParagraphs = [];\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  const cellParas = cell.body.paragraphs;\n  cellParas.load(\"items\");\n  cellParagraphs.push(cellParas);\n}\nawait context.sync();\n\nupdates.forEach((u, i) => {\n  const range = cellParagraphs[i].items[0].getRange();\n  range.insertText(u.text, Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading paragraph\n$d.Paragraphs.Item(1).Range.Text = \"2026-01-30 Friday\"\n\n# Update the division problems in the table (row, column) -> new text\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    \"1,1\" = \"23\u00f75=\"\n    \"1,2\" = \"21\u00f75=\"\n    \"1,3\" = \"27\u00f79=\"\n    \"1,4\" = \"37\u00f78=\"\n    \"1,5\" = \"30\u00f76=\"\n\n    \"5,1\" = \"64\u00f78=\"\n    \"5,2\" = \"22\u00f75=\"\n    \"5,3\" = \"44\u00f73=\"\n    \"5,4\" = \"51\u00f76=\"\n    \"5,5\" = \"89\u00f78=\"\n\n    \"9,1\" = \"75\u00f72=\"\n    \"9,2\" = \"47\u00f73=\"\n    \"9,3\" = \"72\u00f78=\"\n    \"9,4\" = \"42\u00f75=\"\n    \"9,5\" = \"70\u00f75=\"\n\n    \"13,1\" = \"19\u00f76=\"\n    \"13,2\" = \"97\u00f75=\"\n    \"13,3\" = \"61\u00f73=\"\n    \"13,4\" = \"70\u00f79=\"\n    \"13,5\" = \"83\u00f77=\"\n\n    \"17,1\" = \"70\u00f77=\"\n    \"17,2\" = \"99\u00f75=\"\n    \"17,3\" = \"61\u00f77=\"\n    \"17,4\" = \"66\u00f79=\"\n    \"17,5\" = \"99\u00f74=\"\n}\n\nforeach ($row in @(1, 5, 9, 13, 17)) {\n    for ($col = 1; $col -le 5; $col++) {\n        $key = \"$row,$col\"\n        $cell = $t.Cell($row, $col)\n        $cell.Range.Text = $updates[$key]\n    }\n}\n"}
